$wb = $excel.ActiveWorkbook

# Sheets: 1 = Overview, 2 = zh-cn, 3 = de-de
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------
# 1) Update status text from "Ready for handoff" to "In Translation"
#    (appears in the zh-cn / de-de status columns, and their
#    roll-up on the Overview sheet).
# ---------------------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# ---------------------------------------------------------------
# 2) Narrow the status columns to match the shorter text
#    (Overview columns E & F, and column C on zh-cn / de-de).
# ---------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
